## New sea otter index for 2018 (Stephens) built
## Adds 4 new columns to Sheet1: otter_region, date_grass_MM.DD.YY (between the
## site coordinates and the existing survey dates), plus so_duration and
## pop_dens_surv_km2 appended after site_polygon_area_km2. Renames the two
## existing survey-date headers to the new otts1/otts2 naming.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Shift the two existing date columns (D,E -> F,G; old F -> H) right by
# inserting two blank columns at D:E. Column widths/number formats on the
# untouched columns travel with the cells they contain.
$ws.Columns("D:E").Insert()

# --- Header row ------------------------------------------------------------
$ws.Range("D1").Value = "otter_region"
$ws.Range("E1").Value = "date_grass_MM.DD.YY"
$ws.Range("F1").Value = "date_otts1_MM.DD.YY"
$ws.Range("G1").Value = "date_otts2_MM.DD.YY"
$ws.Range("I1").Value = "so_duration"
$ws.Range("J1").Value = "pop_dens_surv_km2"

# --- Data rows ---------------------------------------------------------------
# region, date_grass, so_duration, pop_dens_surv_km2 per site row (2-27)
$rows = @{
   2  = @("low",  43294, 7,  0.154)
   3  = @("mid",  43292, 14, 1.341)
   4  = @("high", 43269, 14, 2.877)
   5  = @("high", 43264, 14, 2.877)
   6  = @("low",  43323, 7,  0.154)
   7  = @("high", 43293, 14, 1.341)
   8  = @("low",  43265, 7,  0.154)
   9  = @("high", 43251, 14, 2.877)
   10 = @("high", 43249, 14, 2.877)
   11 = @("high", 43250, 14, 2.877)
   12 = @("high", 43295, 14, 2.877)
   13 = @("low",  43267, 7,  0.154)
   14 = @("mid",  43279, 14, 1.341)
   15 = @("low",  43266, 7,  0.154)
   16 = @("low",  43267, 7,  0.154)
   17 = @("low",  43266, 7,  0.154)
   18 = @("high", 43239, 14, 2.877)
   19 = @("high", 43281, 14, 2.877)
   20 = @("low",  43291, 7,  1.341)
   21 = @("low",  43322, 7,  0.154)
   22 = @("mid",  43276, 14, 1.341)
   23 = @("mid",  43263, 14, 1.341)
   24 = @("high", 43238, 14, 2.877)
   25 = @("mid",  43268, 14, 1.341)
   26 = @("low",  43278, 7,  0.154)
   27 = @("low",  43265, 7,  0.154)
}

foreach ($r in $rows.Keys) {
   $vals = $rows[$r]
   $ws.Range("D$r").Value = $vals[0]
   $ws.Range("E$r").Value = $vals[1]
   $ws.Range("E$r").NumberFormat = "mm/dd/yy;@"
   $ws.Range("I$r").Value = $vals[2]
   $ws.Range("J$r").Value = $vals[3]
}

# --- Sort range grew by two columns (through H) -----------------------------
$ws.Range("A2:H27").Sort($ws.Range("A2:A27"))
